$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4398919999999999
$ws.Range("I2").Value = 0.719974794695429
$ws.Range("J2").Value = 0.794097205716256
$ws.Range("M2").Value = 2.6796385
$ws.Range("N2").Value = 5.359277000000001
$ws.Range("O2").Value = 0.03934793987875059
$ws.Range("P2").Value = 0.02842274181890429
$ws.Range("Q2").Value = 1.178751539042
$ws.Range("R2").Value = 7.072509234252
$ws.Range("S2").Value = 0.02832952493589154
$ws.Range("T2").Value = 0.02257041985718647

$ws.Range("G3").Value = 0.4398919999999999
$ws.Range("I3").Value = 0.719974794695429
$ws.Range("J3").Value = 0.794097205716256
$ws.Range("M3").Value = 50.102415
$ws.Range("O3").Value = 0.7357062578404556
$ws.Range("P3").Value = 0.7971493203553003
$ws.Range("R3").Value = 198.35686385262
$ws.Range("S3").Value = 0.5296899619448244
$ws.Range("T3").Value = 0.6330140478327565

$ws.Range("G4").Value = 0.4398919999999999
$ws.Range("I4").Value = 0.719974794695429
$ws.Range("J4").Value = 0.794097205716256
$ws.Range("M4").Value = 0.2784063333333333
$ws.Range("N4").Value = 0.8352189999999999
$ws.Range("O4").Value = 0.004088131912518571
$ws.Range("P4").Value = 0.00442955532980352
$ws.Range("Q4").Value = 0.1224687187826666
$ws.Range("R4").Value = 1.102218469044
$ws.Range("S4").Value = 0.00294335193440339
$ws.Range("T4").Value = 0.003517497509962524

$ws.Range("G5").Value = 0.4398919999999999
$ws.Range("I5").Value = 0.719974794695429
$ws.Range("J5").Value = 0.794097205716256
$ws.Range("M5").Value = 13.0677535
$ws.Range("N5").Value = 26.135507
$ws.Range("O5").Value = 0.1918875173156127
$ws.Range("P5").Value = 0.1386087652806835
$ws.Range("Q5").Value = 5.748400222621998
$ws.Range("R5").Value = 34.490401335732
$ws.Range("S5").Value = 0.1381541758839239
$ws.Range("T5").Value = 0.1100688331971712

$ws.Range("G6").Value = 0.4398919999999999
$ws.Range("I6").Value = 0.719974794695429
$ws.Range("J6").Value = 0.794097205716256
$ws.Range("M6").Value = 1.863198333333333
$ws.Range("N6").Value = 5.589594999999999
$ws.Range("O6").Value = 0.02735929342789644
$ws.Range("P6").Value = 0.02964422543511714
$ws.Range("Q6").Value = 0.8196060412466665
$ws.Range("R6").Value = 7.376454371219998
$ws.Range("S6").Value = 0.01969800166876174
$ws.Range("T6").Value = 0.02354039658364929

$ws.Range("G7").Value = 0.4398919999999999
$ws.Range("I7").Value = 0.719974794695429
$ws.Range("J7").Value = 0.794097205716256
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1097013333333333
$ws.Range("N7").Value = 0.329104
$ws.Range("O7").Value = 0.001610859624766094
$ws.Range("P7").Value = 0.001745391780191372
$ws.Range("Q7").Value = 0.04825673892266666
$ws.Range("R7").Value = 0.434310650304
$ws.Range("S7").Value = 0.001159778327624124
$ws.Range("T7").Value = 0.00138601073553009

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.1710905
$ws.Range("H8").Value = 0.342181
$ws.Range("I8").Value = 0.2800252053045709
$ws.Range("J8").Value = 0.2059027942837441
$ws.Range("M8").Value = 2.6796385
$ws.Range("N8").Value = 5.359277000000001
$ws.Range("O8").Value = 0.03934793987875059
$ws.Range("P8").Value = 0.02842274181890429
$ws.Range("Q8").Value = 0.4584606907842501
$ws.Range("R8").Value = 1.833842763137
$ws.Range("S8").Value = 0.01101841494285905
$ws.Range("T8").Value = 0.00585232196171782

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.1710905
$ws.Range("H9").Value = 0.342181
$ws.Range("I9").Value = 0.2800252053045709
$ws.Range("J9").Value = 0.2059027942837441
$ws.Range("M9").Value = 50.102415
$ws.Range("O9").Value = 0.7357062578404556
$ws.Range("P9").Value = 0.7971493203553003
$ws.Range("Q9").Value = 8.5720472335575
$ws.Range("R9").Value = 51.432283401345
$ws.Range("S9").Value = 0.2060162958956312
$ws.Range("T9").Value = 0.1641352725225438

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1710905
$ws.Range("H10").Value = 0.342181
$ws.Range("I10").Value = 0.2800252053045709
$ws.Range("J10").Value = 0.2059027942837441
$ws.Range("M10").Value = 0.2784063333333333
$ws.Range("N10").Value = 0.8352189999999999
$ws.Range("O10").Value = 0.004088131912518571
$ws.Range("P10").Value = 0.00442955532980352
$ws.Range("Q10").Value = 0.04763267877316667
$ws.Range("R10").Value = 0.285796072639
$ws.Range("S10").Value = 0.001144779978115181
$ws.Range("T10").Value = 0.0009120578198409963

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1710905
$ws.Range("H11").Value = 0.342181
$ws.Range("I11").Value = 0.2800252053045709
$ws.Range("J11").Value = 0.2059027942837441
$ws.Range("M11").Value = 13.0677535
$ws.Range("N11").Value = 26.135507
$ws.Range("O11").Value = 0.1918875173156127
$ws.Range("P11").Value = 0.1386087652806835
$ws.Range("Q11").Value = 2.23576848019175
$ws.Range("R11").Value = 8.943073920766999
$ws.Range("S11").Value = 0.05373334143168886
$ws.Range("T11").Value = 0.02853993208351234

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1710905
$ws.Range("H12").Value = 0.342181
$ws.Range("I12").Value = 0.2800252053045709
$ws.Range("J12").Value = 0.2059027942837441
$ws.Range("M12").Value = 1.863198333333333
$ws.Range("N12").Value = 5.589594999999999
$ws.Range("O12").Value = 0.02735929342789644
$ws.Range("P12").Value = 0.02964422543511714
$ws.Range("Q12").Value = 0.3187755344491666
$ws.Range("R12").Value = 1.912653206695
$ws.Range("S12").Value = 0.007661291759134699
$ws.Range("T12").Value = 0.006103828851467858

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1710905
$ws.Range("H13").Value = 0.342181
$ws.Range("I13").Value = 0.2800252053045709
$ws.Range("J13").Value = 0.2059027942837441
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.1097013333333333
$ws.Range("N13").Value = 0.329104
$ws.Range("O13").Value = 0.001610859624766094
$ws.Range("P13").Value = 0.001745391780191372
$ws.Range("Q13").Value = 0.01876885597066667
$ws.Range("R13").Value = 0.112613135824
$ws.Range("S13").Value = 0.0004510812971419694
$ws.Range("T13").Value = 0.0003593810446612819
